# Update the "K" (strikeouts) column (G) for each game row with the
# regenerated values (previously this column held a different stat,
# "Strike#"; the save_data pipeline now regenerates K, std/mean and
# writes s_vals, but only the per-row K values need to change here).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @(1, 1, 0, 0, 1, 1, 0, 1, 1, 7, 1, 6, 7, 1, 4, 5, 7, 4, 8, 1, 4, 5, 3, 9, 0, 5, 2, 5, 6, 10, 6, 6, 2, 4)

$firstRow = 2
$lastRow = $firstRow + $kValues.Length - 1

$values = New-Object 'object[,]' $kValues.Length, 1
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $values[$i, 0] = $kValues[$i]
}

$ws.Range("G$firstRow`:G$lastRow").Value = $values
